$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.792.41'
$ws.Range("E2").Value = '  +1.39%  '

$ws.Range("D3").Value = '3.474.95'
$ws.Range("E3").Value = '  +1.95%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '415.68'
$ws.Range("E5").Value = '  +1.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.08'
$ws.Range("E6").Value = '  +0.81%  '

$ws.Range("E7").Value = '  -0.96%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  -0.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.152'
$ws.Range("E10").Value = '  +6.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.48'
$ws.Range("E11").Value = '  -0.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.78'
$ws.Range("E12").Value = '  +4.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000225'
$ws.Range("E13").Value = '  +1.46%  '

$ws.Range("D14").Value = '4.026.37'
$ws.Range("E14").Value = '  +1.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.56'
$ws.Range("E16").Value = '  -3.11%  '

$ws.Range("D17").Value = '3.468.01'
$ws.Range("E17").Value = '  +1.61%  '

$ws.Range("E18").Value = '  +0.63%  '

$ws.Range("E19").Value = '  -1.30%  '

$ws.Range("D20").Value = '62.767.62'
$ws.Range("E20").Value = '  +1.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '463.33'
$ws.Range("E21").Value = '  +3.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '90.43'
$ws.Range("E22").Value = '  -1.87%  '

$ws.Range("E23").Value = '  +2.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.26'
$ws.Range("E24").Value = '  +0.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.77'
$ws.Range("E25").Value = '  +14.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.32'
$ws.Range("E26").Value = '  +1.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.24'
$ws.Range("E27").Value = '  +0.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.79'
$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.59'
$ws.Range("E29").Value = '  -2.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.12'
$ws.Range("E30").Value = '  +1.12%  '

$ws.Range("E31").Value = '  -3.47%  '

$ws.Range("E32").Value = '  -1.00%  '

$ws.Range("E33").Value = '  -1.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.99'
$ws.Range("E34").Value = '  -4.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.17'
$ws.Range("E36").Value = '  +8.10%  '

$ws.Range("E37").Value = '  -2.80%  '

$ws.Range("E38").Value = '  +0.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.05'
$ws.Range("E39").Value = '  +3.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '148.36'
$ws.Range("E40").Value = '  +3.46%  '

$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.72'
$ws.Range("E41").Value = '  +6.33%  '

$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.36'
$ws.Range("E42").Value = '  -0.59%  '

$ws.Range("E43").Value = '  -0.49%  '

$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.323'
$ws.Range("E44").Value = '  +0.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.42'
$ws.Range("E45").Value = '  +3.30%  '

$ws.Range("E46").Value = '  +2.89%  '

$ws.Range("D47").Value = '0.0₃0586'
$ws.Range("E47").Value = '  +37.43%  '

$ws.Range("E48").Value = '  +10.71%  '

$ws.Range("E49").Value = '  -1.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.34'
$ws.Range("E50").Value = '  +0.35%  '

$ws.Range("E51").Value = '  -3.46%  '
